$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interventions cost and coverage")

# Removed BFP (Breastfeeding promotion, row 2) and CFE (Complementary feeding
# education, originally row 5, now row 4 after the first deletion) from the
# cost & coverage tab.
$ws.Rows.Item(2).Delete() | Out-Null
$ws.Rows.Item(4).Delete() | Out-Null

$ws.Select()
$ws.Range("C7").Select() | Out-Null
